# Applies the "Updated cryptos list" data refresh described by the commit diff.
# Strategy: write new values for the Coin/Link/Price/Volume(1h) cells that changed.
# A handful of Price cells (column D) are purely numeric-looking strings (e.g. "0.997");
# Excel auto-converts those to numbers on assignment, so for just those cells we force
# the Text number format first (shared single style, matches how Excel itself would keep
# the original inlineStr text cells as text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price cells that look like plain numbers: force Text format so they stay strings ---
$numericTextCells = @("D4", "D5", "D6", "D7", "D8", "D15", "D19", "D21", "D23", "D24", "D25", "D27", "D29", "D32", "D33", "D34", "D35", "D36", "D39", "D40", "D42", "D43", "D44", "D45", "D46", "D47", "D49")
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Assign the new cell values (Coin / Link / Price / Volume(1h)) ---
$ws.Range("D4").Value = '0.997'
$ws.Range("D5").Value = '538.05'
$ws.Range("D6").Value = '138.58'
$ws.Range("D7").Value = '0.997'
$ws.Range("D8").Value = '0.566'
$ws.Range("D15").Value = '23.24'
$ws.Range("D19").Value = '11.11'
$ws.Range("D21").Value = '326.44'
$ws.Range("D23").Value = '5.89'
$ws.Range("D24").Value = '65.65'
$ws.Range("D25").Value = '0.424'
$ws.Range("D27").Value = '0.997'
$ws.Range("D29").Value = '6.77'
$ws.Range("D32").Value = '169.61'
$ws.Range("D33").Value = '1.19'
$ws.Range("D34").Value = '0.998'
$ws.Range("D35").Value = '1.47'
$ws.Range("D36").Value = '18.55'
$ws.Range("D39").Value = '36.65'
$ws.Range("D40").Value = '0.826'
$ws.Range("D42").Value = '284.36'
$ws.Range("D43").Value = '5.26'
$ws.Range("D44").Value = '0.999'
$ws.Range("D45").Value = '0.606'
$ws.Range("D46").Value = '130.47'
$ws.Range("D47").Value = '10.85'
$ws.Range("D49").Value = '0.0512'

$ws.Range("D2").Value = '59.236.67'
$ws.Range("E2").Value = '  +0.51%  '
$ws.Range("D3").Value = '2.525.53'
$ws.Range("E3").Value = '  +0.49%  '
$ws.Range("E4").Value = '  -0.37%  '
$ws.Range("E5").Value = '  +1.50%  '
$ws.Range("E6").Value = '  -0.09%  '
$ws.Range("E7").Value = '  -0.27%  '
$ws.Range("E8").Value = '  +0.67%  '
$ws.Range("D9").Value = '2.523.23'
$ws.Range("E9").Value = '  +0.30%  '
$ws.Range("E10").Value = '  +1.43%  '
$ws.Range("E12").Value = '  -1.49%  '
$ws.Range("E13").Value = '  -2.19%  '
$ws.Range("D14").Value = '2.950.82'
$ws.Range("E14").Value = '  -0.19%  '
$ws.Range("E15").Value = '  +1.17%  '
$ws.Range("D16").Value = '59.028.87'
$ws.Range("E16").Value = '  +0.25%  '
$ws.Range("E17").Value = '  +0.19%  '
$ws.Range("D18").Value = '2.519.56'
$ws.Range("E18").Value = '  +0.41%  '
$ws.Range("E19").Value = '  +1.01%  '
$ws.Range("E20").Value = '  +1.36%  '
$ws.Range("E21").Value = '  +1.56%  '
$ws.Range("E22").Value = '  +0.01%  '
$ws.Range("E23").Value = '  +1.49%  '
$ws.Range("E24").Value = '  +5.49%  '
$ws.Range("E25").Value = '  +0.42%  '
$ws.Range("E26").Value = '  +0.04%  '
$ws.Range("E27").Value = '  +0.16%  '
$ws.Range("E28").Value = '  -1.55%  '
$ws.Range("E29").Value = '  +0.56%  '
$ws.Range("D30").Value = '0.0₃0779'
$ws.Range("E30").Value = '  +1.43%  '
$ws.Range("E31").Value = '  +0.30%  '
$ws.Range("E32").Value = '  +4.85%  '
$ws.Range("E33").Value = '  +7.43%  '
$ws.Range("B34").Value = 'USDe'
$ws.Range("C34").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("E34").Value = '  -0.11%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("E35").Value = '  +2.62%  '
$ws.Range("E36").Value = '  +0.76%  '
$ws.Range("E37").Value = '  -1.63%  '
$ws.Range("E38").Value = '  +0.30%  '
$ws.Range("E39").Value = '  -0.81%  '
$ws.Range("E40").Value = '  +3.22%  '
$ws.Range("E41").Value = '  +0.60%  '
$ws.Range("B42").Value = 'Bittensor'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("E42").Value = '  +2.07%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("E43").Value = '  +1.62%  '
$ws.Range("E44").Value = '  +0.07%  '
$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("E45").Value = '  +2.04%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("E46").Value = '  +7.49%  '
$ws.Range("E47").Value = '  -0.13%  '
$ws.Range("E49").Value = '  +0.73%  '
$ws.Range("E50").Value = '  +0.36%  '
$ws.Range("E51").Value = '  +0.34%  '
